$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '45.156.05'
$ws.Cells.Item(2, 5).Value = '  +5.45%  '
$ws.Cells.Item(3, 4).Value = '2.360.73'
$ws.Cells.Item(3, 5).Value = '  +2.00%  '
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '310.61'
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  -0.56%  '
$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '108.12'
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  +0.54%  '
$ws.Cells.Item(7, 5).Value = '  -0.30%  '
$ws.Cells.Item(8, 5).Value = '  -0.06%  '
$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.615'
$ws.Cells.Item(9, 4).Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  +1.12%  '
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '41.14'
$ws.Cells.Item(10, 4).Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  +2.39%  '
$origStyle = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0919'
$ws.Cells.Item(11, 4).Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  +0.40%  '
$origStyle = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '8.45'
$ws.Cells.Item(12, 4).Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  +0.49%  '
$ws.Cells.Item(13, 5).Value = '  +2.01%  '
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.983'
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  -0.96%  '
$ws.Cells.Item(15, 4).Value = '2.715.73'
$ws.Cells.Item(15, 5).Value = '  +1.93%  '
$origStyle = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '15.28'
$ws.Cells.Item(16, 4).Style = $origStyle
$ws.Cells.Item(16, 5).Value = '  +0.04%  '
$ws.Cells.Item(17, 4).Value = '2.360.83'
$ws.Cells.Item(17, 5).Value = '  +2.08%  '
$ws.Cells.Item(18, 4).Value = '45.125.36'
$ws.Cells.Item(18, 5).Value = '  +5.47%  '
$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '14.11'
$ws.Cells.Item(19, 4).Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  +7.62%  '
$ws.Cells.Item(20, 5).Value = '  -2.15%  '
$ws.Cells.Item(21, 5).Value = '  +0.49%  '
$ws.Cells.Item(22, 5).Value = '  -0.74%  '
$ws.Cells.Item(23, 5).Value = '  -0.38%  '
$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '259.91'
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  -2.62%  '
$ws.Cells.Item(25, 5).Value = '  +3.12%  '
$ws.Cells.Item(26, 5).Value = '  -0.38%  '
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.17'
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  +0.85%  '
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.30'
$ws.Cells.Item(28, 4).Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  -5.33%  '
$ws.Cells.Item(29, 5).Value = '  +2.38%  '
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.0967'
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  +10.25%  '
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '22.32'
$ws.Cells.Item(31, 4).Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  -0.82%  '
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '37.78'
$ws.Cells.Item(32, 4).Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  -2.01%  '
$origStyle = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '168.46'
$ws.Cells.Item(33, 4).Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  +1.21%  '
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '2.93'
$ws.Cells.Item(34, 4).Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  +6.98%  '
$ws.Cells.Item(35, 5).Value = '  -0.36%  '
$origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.118'
$ws.Cells.Item(36, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.82'
$ws.Cells.Item(37, 4).Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  +1.85%  '
$ws.Cells.Item(38, 5).Value = '  +5.77%  '
$origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '3.93'
$ws.Cells.Item(39, 4).Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  +6.65%  '
$origStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.0354'
$ws.Cells.Item(40, 4).Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  -1.11%  '
$ws.Cells.Item(41, 5).Value = '  +7.07%  '
$origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '99.04'
$ws.Cells.Item(42, 4).Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  -5.24%  '
$ws.Cells.Item(43, 5).Value = '  -0.88%  '
$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '69.57'
$ws.Cells.Item(44, 4).Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  -1.91%  '
$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '12.85'
$ws.Cells.Item(45, 4).Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  -1.64%  '
$ws.Cells.Item(46, 5).Value = '  +0.20%  '
$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '81.40'
$ws.Cells.Item(47, 4).Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  +5.33%  '
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '111.99'
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  -1.03%  '
$ws.Cells.Item(49, 5).Value = '  +4.81%  '
$ws.Cells.Item(50, 4).Value = '1.682.74'
$ws.Cells.Item(50, 5).Value = '  +1.18%  '
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '9.18'
$ws.Cells.Item(51, 4).Style = $origStyle
$ws.Cells.Item(51, 5).Value = '  +4.24%  '
